$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 316, shifting existing rows 316-346 down to 317-347.
$ws.Rows(316).Insert()

# Populate the newly inserted row 316 with the new weekly record.
$ws.Cells.Item(316, 1).Value = 6
$ws.Cells.Item(316, 2).Value = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Cells.Item(316, 3).Value = "Metropolitana"
$ws.Cells.Item(316, 4).Value = 44461
$ws.Cells.Item(316, 5).Value = 13
$ws.Cells.Item(316, 6).Value = 100112044
$ws.Cells.Item(316, 7).Value = "Perejil"
$ws.Cells.Item(316, 8).Value = "Sin especificar"
$ws.Cells.Item(316, 9).Value = "Primera"
$ws.Cells.Item(316, 10).Value = 150
$ws.Cells.Item(316, 11).Value = 8000
$ws.Cells.Item(316, 12).Value = 9000
$ws.Cells.Item(316, 13).Value = 8400
$ws.Cells.Item(316, 14).Value = "$/docena de atados"
$ws.Cells.Item(316, 15).Value = "Región Metropolitana"
$ws.Cells.Item(316, 16).Value = 2800
$ws.Cells.Item(316, 17).Value = 3
$ws.Cells.Item(316, 18).Value = "Hortaliza"
